# Reorder the "Recorded By" (column G) comma-separated list of recorders so
# that any "System" / "system" entries move to the end of the list, while
# preserving the relative order of the remaining (real) recorder names.
#
# Example: "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#          "System, system, backup@backdoor.com" -> "backup@backdoor.com, System, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count()

# Column G is the 7th column ("Recorded By"); data starts at row 2 (row 1 is the header).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value()

    if ($value -eq $null) {
        continue
    }

    $text = $value.ToString()
    if ($text -eq "") {
        continue
    }

    $parts = $text.Split(",")

    $nonSystem = @()
    $systemOnes = @()
    foreach ($p in $parts) {
        $trimmed = $p.Trim()
        if ($trimmed -eq "System" -or $trimmed -eq "system") {
            $systemOnes += $trimmed
        } else {
            $nonSystem += $trimmed
        }
    }

    # Only rewrite the cell when there's actually a "System" entry that is not
    # already trailing everything else (keeps untouched rows untouched).
    if ($systemOnes.Length -gt 0) {
        $newParts = $nonSystem + $systemOnes
        $newText = [string]::Join(", ", $newParts)

        if ($newText -ne $text) {
            $cell.Value = $newText
        }
    }
}
